# Update the main GSC export data on the "Chart" sheet:
#  - drop the oldest day (2025-11-07), shifting every later day's row up by one
#  - append the newest day (2026-02-04) at the bottom of the series

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest data row (2025-11-07); this shifts all following rows
# up by one, which is exactly the "Invalid"/"Items" shift seen for every
# remaining date.
$ws.Rows("2:2").Delete()

# Append the new trailing day right after the current last data row
# (row 90 once the delete above has collapsed the range by one row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Range("A" + $newRow).Value = "'2026-02-04"
$ws.Range("A" + $newRow).ClearFormats()
$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 28
